$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new D value or $null, new E value)
$changes = @(
    @{ Row = 2; D = "69.013.43"; E = "  +2.80%  " },
    @{ Row = 3; D = "3.742.86"; E = "  +2.39%  " },
    @{ Row = 4; D = "1.00"; E = "  +0.09%  " },
    @{ Row = 5; D = "602.22"; E = "  +1.52%  " },
    @{ Row = 6; D = "168.54"; E = "  +2.22%  " },
    @{ Row = 7; D = "3.742.74"; E = "  +2.51%  " },
    @{ Row = 8; D = $null; E = "  -0.03%  " },
    @{ Row = 9; D = $null; E = "  +2.20%  " },
    @{ Row = 10; D = $null; E = "  +5.67%  " },
    @{ Row = 11; D = $null; E = "  +3.47%  " },
    @{ Row = 12; D = $null; E = "  +0.75%  " },
    @{ Row = 13; D = "38.21"; E = "  +2.66%  " },
    @{ Row = 14; D = "0.0000248"; E = "  +4.31%  " },
    @{ Row = 15; D = "4.369.16"; E = "  +2.23%  " },
    @{ Row = 16; D = "3.750.38"; E = "  +2.46%  " },
    @{ Row = 17; D = "68.972.51"; E = "  +2.64%  " },
    @{ Row = 18; D = $null; E = "  +2.17%  " },
    @{ Row = 19; D = $null; E = "  +0.22%  " },
    @{ Row = 20; D = "17.12"; E = "  -1.06%  " },
    @{ Row = 21; D = "10.88"; E = "  +20.33%  " },
    @{ Row = 22; D = "493.89"; E = "  +0.78%  " },
    @{ Row = 23; D = "0.725"; E = "  +1.69%  " },
    @{ Row = 24; D = "0.0000154"; E = "  +13.87%  " },
    @{ Row = 25; D = "85.26"; E = "  +0.12%  " },
    @{ Row = 26; D = "2.32"; E = "  +1.70%  " },
    @{ Row = 27; D = "12.34"; E = "  +2.13%  " },
    @{ Row = 28; D = "10.41"; E = "  +5.02%  " },
    @{ Row = 29; D = $null; E = "  +0.49%  " },
    @{ Row = 30; D = $null; E = "  +7.73%  " },
    @{ Row = 31; D = "2.97"; E = "  +2.12%  " },
    @{ Row = 32; D = "7.95"; E = "  +4.42%  " },
    @{ Row = 33; D = $null; E = "  +1.49%  " },
    @{ Row = 34; D = "3.889.06"; E = "  +2.45%  " },
    @{ Row = 35; D = $null; E = "  +1.75%  " },
    @{ Row = 36; D = "3.677.30"; E = "  +2.33%  " },
    @{ Row = 38; D = "1.01"; E = "  +3.11%  " },
    @{ Row = 39; D = "5.85"; E = "  +2.22%  " },
    @{ Row = 40; D = $null; E = "  +2.06%  " },
    @{ Row = 41; D = "0.323"; E = "  +0.60%  " },
    @{ Row = 42; D = "3.00"; E = "  +9.19%  " },
    @{ Row = 43; D = "437.50"; E = "  +1.40%  " },
    @{ Row = 44; D = "48.88"; E = "  +0.59%  " },
    @{ Row = 45; D = $null; E = "  +3.36%  " },
    @{ Row = 46; D = "8.44"; E = "  +1.69%  " },
    @{ Row = 48; D = "40.39"; E = "  -0.08%  " },
    @{ Row = 49; D = "141.66"; E = "  -0.28%  " },
    @{ Row = 50; D = "0.0354"; E = "  +2.74%  " },
    @{ Row = 51; D = "2.769.77"; E = "  +1.02%  " }
)

foreach ($chg in $changes) {
    $row = $chg.Row
    if ($null -ne $chg.D) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $chg.D
    }
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $chg.E
}
